$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Tnfsf12"
$ws.Range("C2").Value = "Cd163"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 3.002483333333333
$ws.Range("H2").Value = 9.007449999999999
$ws.Range("I2").Value = 0.1222246438870418
$ws.Range("J2").Value = 0.1222246438870418
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.3881116666666666
$ws.Range("N2").Value = 1.164335
$ws.Range("O2").Value = 0.3557564478344974
$ws.Range("P2").Value = 0.3557564478344974
$ws.Range("Q2").Value = 1.165298810638889
$ws.Range("R2").Value = 10.48768929575
$ws.Range("S2").Value = 0.0434822051470904
$ws.Range("T2").Value = 0.0434822051470904

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Tnfsf12"
$ws.Range("C3").Value = "Cd163"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 3.002483333333333
$ws.Range("H3").Value = 9.007449999999999
$ws.Range("I3").Value = 0.1222246438870418
$ws.Range("J3").Value = 0.1222246438870418
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.3333333333333333
$ws.Range("M3").Value = 0.01985533333333333
$ws.Range("N3").Value = 0.059566
$ws.Range("O3").Value = 0.01820007864721895
$ws.Range("P3").Value = 0.01820007864721895
$ws.Range("Q3").Value = 0.0596153074111111
$ws.Range("R3").Value = 0.5365377667
$ws.Range("S3").Value = 0.002224498131372489
$ws.Range("T3").Value = 0.002224498131372489

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Tnfsf12"
$ws.Range("C4").Value = "Cd163"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 3.002483333333333
$ws.Range("H4").Value = 9.007449999999999
$ws.Range("I4").Value = 0.1222246438870418
$ws.Range("J4").Value = 0.1222246438870418
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.6829806666666666
$ws.Range("N4").Value = 2.048942
$ws.Range("O4").Value = 0.6260434735182836
$ws.Range("P4").Value = 0.6260434735182836
$ws.Range("Q4").Value = 2.050638068655555
$ws.Range("R4").Value = 18.4557426179
$ws.Range("S4").Value = 0.07651794060857887
$ws.Range("T4").Value = 0.07651794060857887

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Tnfsf12"
$ws.Range("C5").Value = "Cd163"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 11.971037
$ws.Range("H5").Value = 35.913111
$ws.Range("I5").Value = 0.4873151894099666
$ws.Range("J5").Value = 0.4873151894099665
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.3881116666666666
$ws.Range("N5").Value = 1.164335
$ws.Range("O5").Value = 0.3557564478344974
$ws.Range("P5").Value = 0.3557564478344974
$ws.Range("Q5").Value = 4.646099121798334
$ws.Range("R5").Value = 41.814892096185
$ws.Range("S5").Value = 0.173365520760285
$ws.Range("T5").Value = 0.173365520760285

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Tnfsf12"
$ws.Range("C6").Value = "Cd163"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 11.971037
$ws.Range("H6").Value = 35.913111
$ws.Range("I6").Value = 0.4873151894099666
$ws.Range("J6").Value = 0.4873151894099665
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.3333333333333333
$ws.Range("M6").Value = 0.01985533333333333
$ws.Range("N6").Value = 0.059566
$ws.Range("O6").Value = 0.01820007864721895
$ws.Range("P6").Value = 0.01820007864721895
$ws.Range("Q6").Value = 0.2376889299806667
$ws.Range("R6").Value = 2.139200369826
$ws.Range("S6").Value = 0.008869174773245789
$ws.Range("T6").Value = 0.008869174773245789

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Tnfsf12"
$ws.Range("C7").Value = "Cd163"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 11.971037
$ws.Range("H7").Value = 35.913111
$ws.Range("I7").Value = 0.4873151894099666
$ws.Range("J7").Value = 0.4873151894099665
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.6829806666666666
$ws.Range("N7").Value = 2.048942
$ws.Range("O7").Value = 0.6260434735182836
$ws.Range("P7").Value = 0.6260434735182836
$ws.Range("Q7").Value = 8.175986830951333
$ws.Range("R7").Value = 73.58388147856199
$ws.Range("S7").Value = 0.3050804938764358
$ws.Range("T7").Value = 0.3050804938764358

# Row 8
$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Tnfsf12"
$ws.Range("C8").Value = "Cd163"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 9.591765666666666
$ws.Range("H8").Value = 28.775297
$ws.Range("I8").Value = 0.3904601667029916
$ws.Range("J8").Value = 0.3904601667029916
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.3881116666666666
$ws.Range("N8").Value = 1.164335
$ws.Range("O8").Value = 0.3557564478344974
$ws.Range("P8").Value = 0.3557564478344974
$ws.Range("Q8").Value = 3.72267615916611
$ws.Range("R8").Value = 33.504085432495
$ws.Range("S8").Value = 0.138908721927122
$ws.Range("T8").Value = 0.138908721927122

# Row 9
$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Tnfsf12"
$ws.Range("C9").Value = "Cd163"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 9.591765666666666
$ws.Range("H9").Value = 28.775297
$ws.Range("I9").Value = 0.3904601667029916
$ws.Range("J9").Value = 0.3904601667029916
$ws.Range("K9").Value = 1
$ws.Range("L9").Value = 0.3333333333333333
$ws.Range("M9").Value = 0.01985533333333333
$ws.Range("N9").Value = 0.059566
$ws.Range("O9").Value = 0.01820007864721895
$ws.Range("P9").Value = 0.01820007864721895
$ws.Range("Q9").Value = 0.1904477045668889
$ws.Range("R9").Value = 1.714029341102
$ws.Range("S9").Value = 0.007106405742600668
$ws.Range("T9").Value = 0.007106405742600668

# Row 10
$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Tnfsf12"
$ws.Range("C10").Value = "Cd163"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 9.591765666666666
$ws.Range("H10").Value = 28.775297
$ws.Range("I10").Value = 0.3904601667029916
$ws.Range("J10").Value = 0.3904601667029916
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 0.6829806666666666
$ws.Range("N10").Value = 2.048942
$ws.Range("O10").Value = 0.6260434735182836
$ws.Range("P10").Value = 0.6260434735182836
$ws.Range("Q10").Value = 6.550990509530442
$ws.Range("R10").Value = 58.95891458577399
$ws.Range("S10").Value = 0.2444450390332689
$ws.Range("T10").Value = 0.2444450390332689
